$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores every value as literal text in the source workbook
# (e.g. "28.404.07", "1.004"). Several of the new prices look like plain numbers
# (e.g. "1.008", "315.02"), so without an explicit text format Excel would silently
# reinterpret them as numeric values. Force those specific cells to Text format
# first so the assigned string is preserved exactly, matching the rest of the sheet.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.290.89'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '1.812.39'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.61%  '
$ws.Range('D5').Value = '315.02'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('D7').Value = '0.5194'
$ws.Range('E7').Value = '  +1.14%  '
$ws.Range('D8').Value = '0.3822'
$ws.Range('E8').Value = '  -2.57%  '
$ws.Range('D9').Value = '0.07926'
$ws.Range('E9').Value = '  +3.63%  '
$ws.Range('D10').Value = '41.90'
$ws.Range('E10').Value = '  +0.76%  '
$ws.Range('D11').Value = '1.102'
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('D12').Value = '6.341'
$ws.Range('E12').Value = '  +1.18%  '
$ws.Range('B13').Value = 'BinanceUSD'
$ws.Range('C13').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D13').Value = '1.008'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').Value = '20.66'
$ws.Range('E14').Value = '  -1.47%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.820.36'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '7.381'
$ws.Range('E16').Value = '  -1.76%  '
$ws.Range('D17').Value = '93.58'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = '0.00001092'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').Value = '0.06613'
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '1.004'
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = '17.43'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').Value = '5.964'
$ws.Range('E22').Value = '  -2.85%  '
$ws.Range('D23').Value = '28.368.86'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').Value = '11.20'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = '2.242'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('D26').Value = '158.02'
$ws.Range('E26').Value = '  +1.22%  '
$ws.Range('D27').Value = '20.55'
$ws.Range('E27').Value = '  -1.09%  '
$ws.Range('D28').Value = '2.027.33'
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('D29').Value = '2.377'
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('D30').Value = '123.87'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('D31').Value = '0.1102'
$ws.Range('E31').Value = '  +1.45%  '
$ws.Range('D32').Value = '1.062'
$ws.Range('E32').Value = '  -4.24%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').Value = '3.679'
$ws.Range('E33').Value = '  +0.59%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '5.606'
$ws.Range('E34').Value = '  -0.77%  '
$ws.Range('D35').Value = '0.07192'
$ws.Range('E35').Value = '  +1.62%  '
$ws.Range('D36').Value = '12.07'
$ws.Range('E36').Value = '  +7.66%  '
$ws.Range('D37').Value = '0.2174'
$ws.Range('E37').Value = '  -1.21%  '
$ws.Range('D38').Value = '0.02317'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').Value = '8.710'
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('D40').Value = '5.041'
$ws.Range('E40').Value = '  -2.31%  '
$ws.Range('D41').Value = '0.6206'
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('D42').Value = '1.174'
$ws.Range('E42').Value = '  +0.39%  '
$ws.Range('D43').Value = '1.388'
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '13.36'
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.6043'
$ws.Range('E45').Value = '  +2.88%  '
$ws.Range('D46').Value = '3.784'
$ws.Range('E46').Value = '  +2.13%  '
$ws.Range('D47').Value = '125.91'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('D48').Value = '1.204'
$ws.Range('E48').Value = '  +0.51%  '
$ws.Range('D49').Value = '1.936'
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('D50').Value = '0.06873'
$ws.Range('E50').Value = '  -0.34%  '
$ws.Range('D51').Value = '73.17'
$ws.Range('E51').Value = '  -0.92%  '
